# Add two new worksheets ("Wyjazdy" and "Przyjazdy") at the end of the
# workbook, each holding a small "Rodzaj / Liczba uczestnikow / Rok" table
# for 2018-2019, and drop the "tabSelected" state from the previously
# active "Granty_przyznane" sheet (Excel moves tabSelected to whichever
# sheet is active when the file is saved - here that becomes "Przyjazdy").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Wyjazdy"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wyjazdy = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wyjazdy.Name = "Wyjazdy"

# Column A: category labels. Fill the six "plain" categories first (rows
# 2-7), then the header (row 1), then the trailing "Inne wyjazdy" category
# (row 8) - this reproduces the original authoring order of the shared
# strings table.
$wyjazdy.Cells.Item(2, 1).Value = "Udział w naukowych imprezach międzynarodowych"
$wyjazdy.Cells.Item(3, 1).Value = "W ramach programu Erasmus+"
$wyjazdy.Cells.Item(4, 1).Value = "Kwerendy"
$wyjazdy.Cells.Item(5, 1).Value = "Wyjazdy na badania naukowe"
$wyjazdy.Cells.Item(6, 1).Value = "Objazdy naukowe"
$wyjazdy.Cells.Item(7, 1).Value = "Ćwiczenia terenowe"

$wyjazdy.Cells.Item(1, 1).Value = "Rodzaj"

$wyjazdy.Cells.Item(8, 1).Value = "Inne wyjazdy"

$wyjazdy.Cells.Item(1, 2).Value = "Liczba uczestników"
$wyjazdy.Cells.Item(1, 3).Value = "Rok"

$wyjazdyRows = @(
    @("Udział w naukowych imprezach międzynarodowych", 1138, 2018),
    @("W ramach programu Erasmus+", 382, 2018),
    @("Kwerendy", 165, 2018),
    @("Wyjazdy na badania naukowe", 171, 2018),
    @("Objazdy naukowe", 133, 2018),
    @("Ćwiczenia terenowe", 50, 2018),
    @("Inne wyjazdy", 330, 2018),
    @("Udział w naukowych imprezach międzynarodowych", 873, 2019),
    @("W ramach programu Erasmus+", 313, 2019),
    @("Kwerendy", 196, 2019),
    @("Wyjazdy na badania naukowe", 164, 2019),
    @("Objazdy naukowe", 163, 2019),
    @("Ćwiczenia terenowe", 102, 2019),
    @("Inne wyjazdy", 400, 2019)
)

for ($i = 0; $i -lt $wyjazdyRows.Count; $i++) {
    $r = $i + 2
    $row = $wyjazdyRows[$i]
    $wyjazdy.Cells.Item($r, 1).Value = $row[0]
    $wyjazdy.Cells.Item($r, 2).Value = $row[1]
    $wyjazdy.Cells.Item($r, 3).Value = $row[2]
}

$wyjazdy.Columns.Item(1).ColumnWidth = 47.28515625
$wyjazdy.Columns.Item(2).ColumnWidth = 18
$wyjazdy.Columns.Item(3).ColumnWidth = 5

$wyjazdy.Range("F6").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Przyjazdy"
# ---------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$przyjazdy = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$przyjazdy.Name = "Przyjazdy"

$przyjazdy.Cells.Item(1, 1).Value = "Rok"
$przyjazdy.Cells.Item(1, 2).Value = "Liczba uczestnikó"
$przyjazdy.Cells.Item(1, 3).Value = "Rok"

$przyjazdyRows = @(
    @("Udział w naukowych imprezach międzynarodowych", 442, 2018),
    @("W ramach programu Erasmus+", 97, 2018),
    @("Kwerendy", 83, 2018),
    @("Wyjazdy na badania naukowe", 40, 2018),
    @("Objazdy naukowe", 30, 2018),
    @("Ćwiczenia terenowe", 28, 2018),
    @("Inne wyjazdy", 51, 2018),
    @("Udział w naukowych imprezach międzynarodowych", 459, 2019),
    @("W ramach programu Erasmus+", 109, 2019),
    @("Kwerendy", 44, 2019),
    @("Wyjazdy na badania naukowe", 41, 2019),
    @("Objazdy naukowe", 35, 2019),
    @("Ćwiczenia terenowe", 35, 2019),
    @("Inne wyjazdy", 83, 2019)
)

for ($i = 0; $i -lt $przyjazdyRows.Count; $i++) {
    $r = $i + 2
    $row = $przyjazdyRows[$i]
    $przyjazdy.Cells.Item($r, 1).Value = $row[0]
    $przyjazdy.Cells.Item($r, 2).Value = $row[1]
    $przyjazdy.Cells.Item($r, 3).Value = $row[2]
}

$przyjazdy.Columns.Item(1).ColumnWidth = 47.28515625
$przyjazdy.Columns.Item(2).ColumnWidth = 16.28515625
$przyjazdy.Columns.Item(3).ColumnWidth = 5

$przyjazdy.Range("E12").Select() | Out-Null
